# Fix typo on schematic
#   - "5-days of 15% PEG 3350 in drinking water (N-21)" -> "...(N=21)"
#   - ". (N = 9)"                                        -> ". (N=9)"
#   - "5-day PEG 3350 + 10-day recovery (N = 12)"         -> split into
#       "5-day PEG 3350 + 10-day recovery " + "(N=12" + ")"
#
# All three edits live in the same legend textbox on slide 1. Locate it by
# content (rather than a hard-coded index) so the script is resilient to
# shape ordering.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$sh = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $cand = $s.Shapes.Item($i)
    if ($cand.HasTextFrame -and $cand.TextFrame.HasText) {
        if ($cand.TextFrame.TextRange.Text -like "*N-21*") {
            $sh = $cand
        }
    }
}

$tr = $sh.TextFrame.TextRange

# --- Paragraph 2: "5-days of 15% PEG 3350 in drinking water (N-21)" ---
$para2 = $tr.Paragraphs(2, 1)
$run2a = $para2.Characters(1, 47)
$run2a.Text = "5-days of 15% PEG 3350 in drinking water (N=21)"

# --- Paragraph 3: trailing run "...Clind. (N = 9)" ---
$para3 = $tr.Paragraphs(3, 1)
$run3c = $para3.Characters(23, 9)
$run3c.Text = ". (N=9)"

# --- Paragraph 4: "5-day PEG 3350 + 10-day recovery (N = 12)" ---
# First split off the trailing "(N = 12)" into its own run and fix the spacing.
$para4 = $tr.Paragraphs(4, 1)
$run4b = $para4.Characters(34, 8)
$run4b.Text = "(N=12)"

# Then split that run again so the closing parenthesis is its own run.
$run4c = $tr.Paragraphs(4, 1).Characters(39, 1)
$run4c.Text = ")"
